$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row cells: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
$ws.Range("A1").Value2 = "Segmentname_FV2310"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2310"
$ws.Range("C1").Value2 = "Segment_FV2310"
$ws.Range("D1").Value2 = "Datenelement_FV2310"
$ws.Range("E1").Value2 = "Segment ID_FV2310"
$ws.Range("F1").Value2 = "Code_FV2310"
$ws.Range("G1").Value2 = "Qualifier_FV2310"
$ws.Range("H1").Value2 = "Beschreibung_FV2310"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value2 = "Bedingung_FV2310"

$ws.Range("L1").Value2 = "Segmentname_FV2404"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2404"
$ws.Range("N1").Value2 = "Segment_FV2404"
$ws.Range("O1").Value2 = "Datenelement_FV2404"
$ws.Range("P1").Value2 = "Segment ID_FV2404"
$ws.Range("Q1").Value2 = "Code_FV2404"
$ws.Range("R1").Value2 = "Qualifier_FV2404"
$ws.Range("S1").Value2 = "Beschreibung_FV2404"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value2 = "Bedingung_FV2404"

# 2) Freeze the header row (split after row 1) and set the active pane/selection
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the data range into an Excel Table (ListObject) without a table style
$rng = $ws.Range("A1:U65")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
